$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header string updates (shared strings used by the date-stamped columns)
$ws.Range("F1").Value = "AC_20200816"
$ws.Range("G1").Value = "AN_20200816"
$ws.Range("N1").Value = "AN_COTA_20200816"

# Data updates (C,D,E,F,G,N columns for rows 2-42)
$ws.Range("C2").Value = -26.2
$ws.Range("D2").Value = -33
$ws.Range("E2").Value = -24.2
$ws.Range("F2").Value = 45
$ws.Range("G2").Value = -27.2
$ws.Range("N2").Value = -10.78
$ws.Range("C3").Value = -31.1
$ws.Range("D3").Value = -43
$ws.Range("E3").Value = -8.1
$ws.Range("F3").Value = 63.5
$ws.Range("G3").Value = 44.6
$ws.Range("N3").Value = -46.38
$ws.Range("C4").Value = -36.9
$ws.Range("D4").Value = -41.2
$ws.Range("E4").Value = -22
$ws.Range("F4").Value = 46.9
$ws.Range("G4").Value = -40.4
$ws.Range("C5").Value = -10.4
$ws.Range("D5").Value = -6.1
$ws.Range("E5").Value = -0.3
$ws.Range("F5").Value = 84.3
$ws.Range("G5").Value = 33.6
$ws.Range("C6").Value = -24.9
$ws.Range("D6").Value = -32.3
$ws.Range("E6").Value = -32.3
$ws.Range("F6").Value = 83.59999999999999
$ws.Range("G6").Value = -6.1
$ws.Range("N6").Value = -12.87
$ws.Range("C7").Value = -20.8
$ws.Range("D7").Value = -27.4
$ws.Range("E7").Value = -20.6
$ws.Range("F7").Value = 90.7
$ws.Range("G7").Value = -5.7
$ws.Range("C8").Value = -54.3
$ws.Range("D8").Value = -49.8
$ws.Range("E8").Value = -20.5
$ws.Range("F8").Value = 73.3
$ws.Range("G8").Value = -12.9
$ws.Range("C9").Value = -29.1
$ws.Range("D9").Value = -32
$ws.Range("E9").Value = -32.4
$ws.Range("F9").Value = 66.7
$ws.Range("G9").Value = -26
$ws.Range("C10").Value = -50
$ws.Range("D10").Value = -58.5
$ws.Range("E10").Value = -38.1
$ws.Range("F10").Value = 65.2
$ws.Range("G10").Value = -15.7
$ws.Range("N10").Value = -28.1
$ws.Range("C11").Value = -46.1
$ws.Range("D11").Value = -52.8
$ws.Range("E11").Value = -29
$ws.Range("F11").Value = 54
$ws.Range("G11").Value = -36.3
$ws.Range("N11").Value = -15.87
$ws.Range("C12").Value = -46.3
$ws.Range("D12").Value = -53.6
$ws.Range("E12").Value = -29.4
$ws.Range("F12").Value = 54.6
$ws.Range("G12").Value = -37.1
$ws.Range("N12").Value = -15.87
$ws.Range("C13").Value = -25.2
$ws.Range("D13").Value = -31.2
$ws.Range("E13").Value = -31
$ws.Range("F13").Value = 62.5
$ws.Range("G13").Value = -36
$ws.Range("C14").Value = -51.5
$ws.Range("D14").Value = -58.2
$ws.Range("E14").Value = -25.7
$ws.Range("F14").Value = 46.9
$ws.Range("G14").Value = -35.1
$ws.Range("N14").Value = -11.59
$ws.Range("C15").Value = -21.9
$ws.Range("D15").Value = -21.6
$ws.Range("E15").Value = 10.4
$ws.Range("F15").Value = 71.40000000000001
$ws.Range("G15").Value = -26
$ws.Range("N15").Value = -14.64
$ws.Range("C16").Value = -17.3
$ws.Range("D16").Value = -9.5
$ws.Range("E16").Value = 21.1
$ws.Range("F16").Value = 89
$ws.Range("G16").Value = -17.2
$ws.Range("N16").Value = 1.31
$ws.Range("C17").Value = -48.6
$ws.Range("D17").Value = -46.3
$ws.Range("E17").Value = -19
$ws.Range("F17").Value = 43.5
$ws.Range("G17").Value = -51.1
$ws.Range("C18").Value = -39.6
$ws.Range("D18").Value = -42.5
$ws.Range("E18").Value = -34.4
$ws.Range("F18").Value = 57
$ws.Range("G18").Value = -34.2
$ws.Range("N18").Value = -10.78
$ws.Range("C19").Value = -29.1
$ws.Range("D19").Value = -22
$ws.Range("E19").Value = -5.8
$ws.Range("F19").Value = 72.40000000000001
$ws.Range("G19").Value = -25.4
$ws.Range("C20").Value = -40.7
$ws.Range("D20").Value = -44.6
$ws.Range("E20").Value = -29
$ws.Range("F20").Value = 47.2
$ws.Range("G20").Value = -44.2
$ws.Range("C21").Value = -17.9
$ws.Range("D21").Value = -6.9
$ws.Range("E21").Value = 25.6
$ws.Range("F21").Value = 87.59999999999999
$ws.Range("G21").Value = -21
$ws.Range("C22").Value = -29.2
$ws.Range("D22").Value = -20
$ws.Range("E22").Value = 14.6
$ws.Range("F22").Value = 119.4
$ws.Range("G22").Value = 19.5
$ws.Range("C23").Value = -26.7
$ws.Range("D23").Value = -31.1
$ws.Range("E23").Value = -29.2
$ws.Range("F23").Value = 76.5
$ws.Range("G23").Value = -20.5
$ws.Range("C24").Value = -29.6
$ws.Range("D24").Value = -34
$ws.Range("E24").Value = -26.6
$ws.Range("F24").Value = 100
$ws.Range("G24").Value = 33
$ws.Range("C25").Value = -45.2
$ws.Range("D25").Value = -56.8
$ws.Range("E25").Value = -42.3
$ws.Range("F25").Value = 62.6
$ws.Range("G25").Value = -26.4
$ws.Range("C26").Value = -31.8
$ws.Range("D26").Value = -40.5
$ws.Range("E26").Value = -23.5
$ws.Range("F26").Value = 39.7
$ws.Range("G26").Value = -36
$ws.Range("C27").Value = -52.8
$ws.Range("D27").Value = -53.8
$ws.Range("E27").Value = -28.6
$ws.Range("F27").Value = 55.7
$ws.Range("G27").Value = -40
$ws.Range("C28").Value = -14.4
$ws.Range("D28").Value = -19.9
$ws.Range("E28").Value = -11.9
$ws.Range("F28").Value = 74.90000000000001
$ws.Range("G28").Value = -3.3
$ws.Range("C29").Value = -46.2
$ws.Range("D29").Value = -51.5
$ws.Range("E29").Value = -15.2
$ws.Range("F29").Value = 64.90000000000001
$ws.Range("G29").Value = -34
$ws.Range("N29").Value = -34.39
$ws.Range("C30").Value = -7.3
$ws.Range("D30").Value = -1
$ws.Range("E30").Value = 31.5
$ws.Range("F30").Value = 89.40000000000001
$ws.Range("G30").Value = -12.1
$ws.Range("C31").Value = -18.1
$ws.Range("D31").Value = -22.6
$ws.Range("E31").Value = -19
$ws.Range("F31").Value = 71.7
$ws.Range("G31").Value = -16.1
$ws.Range("C32").Value = -14
$ws.Range("D32").Value = -19.7
$ws.Range("E32").Value = -15.2
$ws.Range("F32").Value = 71.2
$ws.Range("G32").Value = -7.5
$ws.Range("C33").Value = -32.7
$ws.Range("D33").Value = -32.7
$ws.Range("E33").Value = -20.8
$ws.Range("F33").Value = 61.7
$ws.Range("G33").Value = -28.6
$ws.Range("C34").Value = -50.8
$ws.Range("D34").Value = -57.1
$ws.Range("E34").Value = -23.7
$ws.Range("F34").Value = 47.2
$ws.Range("G34").Value = -31.8
$ws.Range("N34").Value = -11.59
$ws.Range("C35").Value = -54.8
$ws.Range("D35").Value = -54.9
$ws.Range("E35").Value = -29.7
$ws.Range("F35").Value = 58
$ws.Range("G35").Value = -42.5
$ws.Range("C36").Value = -59.4
$ws.Range("D36").Value = -61.9
$ws.Range("E36").Value = -30.6
$ws.Range("F36").Value = 31.2
$ws.Range("G36").Value = -70.09999999999999
$ws.Range("N36").Value = -70.86
$ws.Range("C37").Value = -11.3
$ws.Range("D37").Value = -12
$ws.Range("E37").Value = 11.7
$ws.Range("F37").Value = 78.09999999999999
$ws.Range("G37").Value = -15
$ws.Range("C38").Value = -42.9
$ws.Range("D38").Value = -53.7
$ws.Range("E38").Value = -30.2
$ws.Range("F38").Value = 63.3
$ws.Range("G38").Value = -19.8
$ws.Range("N38").Value = -28.1
$ws.Range("C39").Value = -5.6
$ws.Range("D39").Value = 3.1
$ws.Range("E39").Value = 36.6
$ws.Range("F39").Value = 92.2
$ws.Range("G39").Value = -13.8
$ws.Range("C40").Value = -13.5
$ws.Range("D40").Value = -18.3
$ws.Range("E40").Value = -0.6
$ws.Range("F40").Value = 74.59999999999999
$ws.Range("G40").Value = -12.7
$ws.Range("C41").Value = -49.4
$ws.Range("D41").Value = -53.5
$ws.Range("E41").Value = -27.3
$ws.Range("F41").Value = 51.9
$ws.Range("G41").Value = -35.5
$ws.Range("N41").Value = -11.59
$ws.Range("C42").Value = -26.8
$ws.Range("D42").Value = -36.3
$ws.Range("E42").Value = -34.2
$ws.Range("F42").Value = 101.3
$ws.Range("G42").Value = 3.9
$ws.Range("N42").Value = -12.87
